$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newText = @"
questions = [
    {
        "title": "Your organization has a mailbox server on Microsoft Exchange Server 2019 with around 1,700 client mailboxes. As a system administrator, you need to ensure that all users receive a notification when they try to send an email message to a specific distribution group. Which cmdlet should you use to set up this notification?",
        "ques_type": 2,
        "options": [
            "New-TransportRule",
            "Get-MailboxDistributionGroup",
            "Apply-MailboxServerParameter",
            "New-NotificationPolicy"
        ],
        "score": "New-TransportRule"
    },
    {
        "title": "You want to copy the end-user data for a user named Alain from one Exchange Server to a mailbox named Ex001 on another Exchange Server. Both servers are in the same organization.Which command should you run?",
        "ques_type": 2,
        "options": [
            "Restore-Mailbox -Identity Alain -RecoveryMailboxDatabase Ex0001",
            "Restore-UserMailbox -Identity Alain -RecoveryDatabase Ex001",
            "Restore-Mailbox -Identity Alain -RecoveryDatabase Ex001",
            "Restore-Mailbox -Identity Alain -RecoveryMailboxDatabase Ex001"
        ],
        "score": "Restore-Mailbox -Identity Alain -RecoveryDatabase Ex001"
    },
    {
        "title": "You are a system administrator. Users start complaining that sent emails are not reaching recipients. You check the mail server and see several emails in the queue. At the top of the list is an email with attachments totaling 12GB (the attachment size limit at your organization is 15GB). You delete this email to free the queue and allow users to work. What immediate action should you take to ensure that incidents of this type do not happen again?",
        "ques_type": 2,
        "options": [
            "Delete queues on the mail server.",
            "Apply a filter to delete by default all emails with more than 1GB worth of attachments in the queue.",
            "Limit the user mailbox storage quota to 12GB.",
            "Reduce the maximum size of email attachments to 25 MB."
        ],
        "score": "Reduce the maximum size of email attachments to 25 MB."
    },
    {
        "title": "True or false: If transport services are down, users from within the same MS Exchange network will be able to send emails to and receive emails from one another.",
        "ques_type": 11,
        "options": [
            "true",
            "false"
        ],
        "score": "False"
    }
]
"@

# Remove the second row entirely (old shared-string cell lived at A2)
$ws.Range("A2").ClearContents()

# Reset A1 back to the default (unstyled) look before writing the new text,
# so no leftover bold/bordered style index gets attached to the cell
$ws.Range("A1").ClearFormats()
$ws.Range("A1").Value = $newText

# Drop the auto row-height Excel applies for wrapped multi-line text so the
# row goes back to using the sheet's default height (no explicit ht/customHeight)
$ws.Rows(1).AutoFit()
